$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 48

$ws.Cells.Item($row, 1).Value = 47
$ws.Cells.Item($row, 2).Value = "Total Memory: 15.86 GB, Used Memory: 7.47 GB, Total Disk Space: 237.84 GB"
$ws.Cells.Item($row, 3).Value = "Sovan.Souern"
$ws.Cells.Item($row, 4).Value = "1L0N1W2"
$ws.Cells.Item($row, 5).Value = "AMD64"
$ws.Cells.Item($row, 6).Value = "Windows"
$ws.Cells.Item($row, 7).Value = "PNCL114"
$ws.Cells.Item($row, 8).Value = "AT/AT COMPATIBLE"

$rowRange = $ws.Range("A" + $row + ":H" + $row)
$rowRange.HorizontalAlignment = -4108
